# "fix dashboard utiblil choeun" - the meter-import dashboard's header for
# column D was mislabeled "customer" when the sheet actually carries the
# property name, so relabel it to "property".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D1").Value = "property"

# keep the active-cell/selection in sync with the author's saved view
$ws.Range("H9").Select()
